# Generate Report for Handoff
#
# A new handoff/handback round has run for "xinjiang":
#   - the markdown source file was regenerated with a new GUID
#     (dfd096ee-948e-4db1-9687-dda1ece87f01.md -> 1214c0cf-72f0-40f9-ad5a-bacb0e69c537.md)
#   - a second, not-localized markdown file shows up
#     (ffff724010a3-8cf2-48a9-aaa3-721a6a52bc3d.md) sitting between the
#     handoff-ready file and the .localization-config entry
#   - the zh-cn / de-de xlf handoff packages were regenerated with a new
#     content hash and new handoff timestamps
#
# This pushes ".localization-config" from row 3 down to row 4 on every
# sheet, and inserts a duplicate-shaped "Ready for handoff" row (backed by
# the new ignored-file GUID) as the new row 3.

$wb = $excel.ActiveWorkbook

$oldGuid = "dfd096ee-948e-4db1-9687-dda1ece87f01"
$newGuid = "1214c0cf-72f0-40f9-ad5a-bacb0e69c537"
$ignoredGuid = "ffff724010a3-8cf2-48a9-aaa3-721a6a52bc3d"
$oldHash = "38819295ffc34122546952c586aec66288bd7449"
$newHash = "9202164902571d7a033129dc8e57e3d8a2cd1b5d"

$newMdName = "$newGuid.md"
$ignoredMdName = "$ignoredGuid.md"
$configName = ".localization-config"

$zhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$deXlfName = "$newGuid.$newHash.de-de.xlf"

$zhTime = "2016-03-07 02:39:18"
$deTime = "2016-03-07 02:39:29"
$epoch = "0001-01-01 00:00:00"

$readyForHandoff = "Ready for handoff"
$notLocalized = "Not to be localized"
$includeTxt = "Include"
$ignoredTxt = "Ignored"

# Hyperlink font used throughout the workbook for linked filenames
# (underline + custom blue FF6495ED) -- matches the "HyperLink" cell style
# already used by every other linked cell in the sheet.
$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) packed as a VBA-style BGR long

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2 ("...md" / Ready for handoff / Ready for handoff) keeps its shape,
# only the backing filename changes.
Style-AsHyperlink($wsOverview.Range("A2"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/e2e/$newMdName", "", "", $newMdName) | Out-Null
$wsOverview.Range("B2").Value = $readyForHandoff
$wsOverview.Range("C2").Value = $readyForHandoff

# Row 3 becomes the new "ignored" markdown file, shaped just like row 2.
Style-AsHyperlink($wsOverview.Range("A3"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/e2e/$ignoredMdName", "", "", $ignoredMdName) | Out-Null
$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff

# Row 4 (new): ".localization-config" moves here.
Style-AsHyperlink($wsOverview.Range("A4"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/$configName", "", "", $configName) | Out-Null
$wsOverview.Range("B4").Value = $notLocalized
$wsOverview.Range("C4").Value = $notLocalized

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2
Style-AsHyperlink($wsZh.Range("A2"))
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/e2e/$newMdName", "", "", $newMdName) | Out-Null
$wsZh.Range("B2").Value = $readyForHandoff
Style-AsHyperlink($wsZh.Range("C2"))
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b77f908d67eac1de7a22e327421faaf363a71f34/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName", "", "", $zhXlfName) | Out-Null
Style-AsDate($wsZh.Range("D2"))
$wsZh.Range("D2").Value = $zhTime
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = $includeTxt

# Row 3 (new shape: duplicate of row 2, backed by the ignored-file guid)
Style-AsHyperlink($wsZh.Range("A3"))
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/e2e/$ignoredMdName", "", "", $ignoredMdName) | Out-Null
$wsZh.Range("B3").Value = $readyForHandoff
Style-AsHyperlink($wsZh.Range("C3"))
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b77f908d67eac1de7a22e327421faaf363a71f34/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName", "", "", $zhXlfName) | Out-Null
Style-AsDate($wsZh.Range("D3"))
$wsZh.Range("D3").Value = $zhTime
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = $includeTxt

# Row 4 (new): ".localization-config"
Style-AsHyperlink($wsZh.Range("A4"))
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/$configName", "", "", $configName) | Out-Null
$wsZh.Range("B4").Value = $notLocalized
Style-AsDate($wsZh.Range("D4"))
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = $ignoredTxt

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2
Style-AsHyperlink($wsDe.Range("A2"))
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/e2e/$newMdName", "", "", $newMdName) | Out-Null
$wsDe.Range("B2").Value = $readyForHandoff
Style-AsHyperlink($wsDe.Range("C2"))
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5aadc6b478a7be8d50e48b4d7170a9bd71efc2f1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName", "", "", $deXlfName) | Out-Null
Style-AsDate($wsDe.Range("D2"))
$wsDe.Range("D2").Value = $deTime
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = $includeTxt

# Row 3 (new shape: duplicate of row 2, backed by the ignored-file guid)
Style-AsHyperlink($wsDe.Range("A3"))
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/e2e/$ignoredMdName", "", "", $ignoredMdName) | Out-Null
$wsDe.Range("B3").Value = $readyForHandoff
Style-AsHyperlink($wsDe.Range("C3"))
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5aadc6b478a7be8d50e48b4d7170a9bd71efc2f1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName", "", "", $deXlfName) | Out-Null
Style-AsDate($wsDe.Range("D3"))
$wsDe.Range("D3").Value = $deTime
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = $includeTxt

# Row 4 (new): ".localization-config"
Style-AsHyperlink($wsDe.Range("A4"))
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/746c351d11e54d3aa0697cfc05f3a1c02be85e29/$configName", "", "", $configName) | Out-Null
$wsDe.Range("B4").Value = $notLocalized
Style-AsDate($wsDe.Range("D4"))
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = $ignoredTxt

Write-Host "Report regenerated for handoff."
